# Update the asset inventory sheet: refresh existing rows 2-9 and append new
# rows 10-19 for the "inactive asset and report log" feature.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(2, "IT653459", "Chair", "dorm_amenities", "h17", "'12190025", 45077),
  @(3, "IT653460", "Table", "dorm_amenities", "h17", "'12190029", 45078),
  @(4, "IT653461", "scoreboard", "dorm_amenities", "SSO office", "'12190026", 45079),
  @(5, "IT653462", "filter", "dorm_amenities", "h17", "'12190040", 45080),
  @(6, "IT65329324", "table", "dorm_amenities", "SSO office", "'12190023", 45074),
  @(7, "IT65329325", "table", "dorm_amenities", "h17", "'12190023", 45075),
  @(8, "IT65329326", "chair", "dorm_amenities", "h18", "'12190024", 45076),
  @(9, "IT65329327", "Chair", "dorm_amenities", "h17", "'12190025", 45077),
  @(10, "IT657348", "Chair", "dorm_amenities", "h17", "'12190025", 45077),
  @(11, "IT657347", "chair", "dorm_amenities", "Y18", "'12190024", 45076),
  @(12, "IT657346", "table", "dorm_amenities", "h17", "'12190023", 45075),
  @(13, "IT65329328", "Table", "dorm_amenities", "h17", "'12190029", 45078),
  @(14, "IT657349", "Table", "dorm_amenities", "h17", "'12190029", 45078),
  @(15, "IT657350", "scoreboard", "dorm_amenities", "SSO office", "'12190026", 45079),
  @(16, "IT65329329", "scoreboard", "dorm_amenities", "SSO office", "'12190026", 45079),
  @(17, "IT657345", "Bed type A", "dorm_amenities", "SSO office", "'12190023", 45074),
  @(18, "IT653456", "table", "dorm_amenities", "SSO office", "'12190023", 45074),
  @(19, "IT653457", "table", "dorm_amenities", "h17", "'12190023", 45075)
)

foreach ($row in $rows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
  $ws.Cells.Item($r, 6).Value = $row[6]
  $ws.Cells.Item($r, 6).NumberFormat = "yyyy-mm-dd"
}
